$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.796.43'
$ws.Range("D3").Value = '1.635.03'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("E8").Value = '  +0.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.64'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0789'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.27'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").Value = '1.636.77'
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").Value = '1.860.23'
$ws.Range("E15").Value = '  -0.18%  '
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.90'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").Value = '25.801.80'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("E20").Value = '  +1.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '194.44'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.28'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.00%  '
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("E25").Value = '  +2.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.84'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.04%  '
$ws.Range("E27").Value = '  -0.51%  '
$ws.Range("E28").Value = '  +0.68%  '
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("E30").Value = '  -0.09%  '
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.35'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.20%  '
$ws.Range("E33").Value = '  +0.26%  '
$ws.Range("E34").Value = '  +1.05%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("D37").Value = '1.130.70'
$ws.Range("E38").Value = '  -1.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.548'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.41%  '
$ws.Range("E40").Value = '  -0.23%  '
$ws.Range("E41").Value = '  +0.50%  '
$ws.Range("E42").Value = '  +2.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.39'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.807'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.93%  '
$ws.Range("D45").Value = '1.769.70'
$ws.Range("E45").Value = '  -0.22%  '
$ws.Range("D46").Value = '0.0₆0112'
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.24'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.46%  '
$ws.Range("E48").Value = '  -2.29%  '
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.57'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.33'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.54%  '
